# ---------------------------------------------------------------------------
# Gastos2020-2021.xlsx — "atualizado planilha de gastor"
#   1. Janeiro: move two mis-typed "Não gasto" values (C14,C15) into the
#      "Valor" column (B14,B15) -- the SUM formulas recompute automatically.
#   2. Add a new "Fevereiro" sheet (Compras 02/2021) mirroring the layout of
#      the existing month sheets, with its own products/prices and totals.
#   3. Move the workbook's active tab to the new Fevereiro sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Janeiro corrections
# ---------------------------------------------------------------------------
$jan = $wb.Worksheets.Item("Janeiro")

# C14 (124) really belongs in column B ("Valor"), same for C15 (95)
$jan.Range("B14").Value = $jan.Range("C14").Value2
$jan.Range("C14").Clear()

$jan.Range("B15").Value = $jan.Range("C15").Value2
$jan.Range("C15").Clear()

# Janeiro is no longer the tab in focus; selection settles on F10
$jan.Range("F10").Select()

# ---------------------------------------------------------------------------
# 2) New "Fevereiro" sheet, inserted after Janeiro (becomes the last sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$feb = $wb.Worksheets.Add($null, $lastSheet)
$feb.Name = "Fevereiro"

# column widths (approximate the source workbook's custom widths)
$feb.Columns("A").ColumnWidth = 21.93
$feb.Columns("B").ColumnWidth = 12.75
$feb.Columns("C").ColumnWidth = 16.93
$feb.Columns("D").ColumnWidth = 15.93
$feb.Columns("E").ColumnWidth = 12.59
$feb.Columns("F").ColumnWidth = 14.43
$feb.Columns("G").ColumnWidth = 10.59

$moedaNovo = '_-[$R$-416]\ * #,##0.00_-;\-[$R$-416]\ * #,##0.00_-;_-[$R$-416]\ * "-"??_-;_-@_-'

# --- header row -------------------------------------------------------
$feb.Range("A1").Value = "Produto"
$feb.Range("D1").Value = "Observação"

$feb.Range("B1").Value = "Valor"
$feb.Range("C1").Value = "Não gasto"
$feb.Range("F1").Value = "Total Gasto"
$hdr = $feb.Range("B1:C1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.NumberFormat = $moedaNovo
$feb.Range("F1").Font.Bold = $true
$feb.Range("F1").HorizontalAlignment = -4108
$feb.Range("F1").NumberFormat = $moedaNovo

$feb.Range("A1:D1").Font.Size = 9
$feb.Range("A1").HorizontalAlignment = -4108
$feb.Range("D1").HorizontalAlignment = -4108
$feb.Range("F1").Font.Size = 9

# --- products / values --------------------------------------------------
$feb.Range("A2").Value = "Mochila"
$feb.Range("C2").Value = 65.69
$feb.Range("C2").Style = "Moeda"

$feb.Range("A3").Value = "Desodorante kaik urbe"
$feb.Range("B3").Style = "Moeda"
$feb.Range("C3").Value = 44
$feb.Range("C3").Style = "Moeda"

$feb.Range("A4").Value = "Cartão SD"
$feb.Range("C4").Value = 109
$feb.Range("C4").Font.Bold = $false
$feb.Range("C4").NumberFormat = $moedaNovo

$feb.Range("A5").Value = "Capa flip iPhone 5s"
$feb.Range("C5").Value = 32.89
$feb.Range("C5").NumberFormat = $moedaNovo

$feb.Range("A6").Value = "Mochila"
$feb.Range("C6").Value = 64.34
$feb.Range("C6").NumberFormat = $moedaNovo

$feb.Range("A7").Value = "Violino elétrico"
$feb.Range("C7").Value = 640
$feb.Range("C7").NumberFormat = $moedaNovo

$feb.Range("A8").Value = "Organizador de cabos"
$feb.Range("C8").Value = 32
$feb.Range("C8").NumberFormat = $moedaNovo

$feb.Range("A9").Value = "IPTU"
$feb.Range("C9").Value = 631.52
$feb.Range("C9").NumberFormat = $moedaNovo

# --- "Compras 02/2021" label box (F3:G7) ---------------------------------
$feb.Range("F3:G7").Merge()
$feb.Range("F3").Value = "Compras 02/2021"
$feb.Range("F3:G7").Font.Bold = $true
$feb.Range("F3:G7").Font.Size = 16
$feb.Range("F3:G7").HorizontalAlignment = -4108
$feb.Range("F3:G7").VerticalAlignment = -4108

# --- totals ---------------------------------------------------------------
$feb.Range("F9").Value = "À Gastar"
$feb.Range("F9").Font.Size = 9
$feb.Range("F9").HorizontalAlignment = -4108

$feb.Range("F2").Formula = "=SUM(B:B)"
$feb.Range("F2").NumberFormat = $moedaNovo

$feb.Range("F10").Formula = "=SUM(C:C)"
$feb.Range("F10").NumberFormat = $moedaNovo

$feb.Range("C10").Select()
$feb.Activate()
